$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '49.970.20'
$ws.Range("E2").Value = '  +3.80%  '
$ws.Range("D3").Value = '2.646.06'
$ws.Range("E3").Value = '  +5.61%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '110.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.531'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +3.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.76'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.61'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0822'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("E13").Value = '  +0.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.79%  '
$ws.Range("D15").Value = '3.049.81'
$ws.Range("E15").Value = '  +5.29%  '
$ws.Range("D16").Value = '2.608.74'
$ws.Range("E16").Value = '  +4.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.881'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.40%  '
$ws.Range("D18").Value = '49.924.91'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.12%  '
$ws.Range("E20").Value = '  +10.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.82%  '
$ws.Range("D22").Value = '0.0₃0966'
$ws.Range("E22").Value = '  +2.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '282.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.66%  '
$ws.Range("E28").Value = '  +6.97%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.75'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.11%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.87%  '
$ws.Range("E31").Value = '  +2.96%  '
$ws.Range("E32").Value = '  +0.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0798'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.99%  '
$ws.Range("E37").Value = '  +6.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.76'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.61%  '
$ws.Range("E39").Value = '  +8.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '124.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.91%  '
$ws.Range("E41").Value = '  +1.52%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.94%  '
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.23'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.00%  '
$ws.Range("E44").Value = '  +4.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.41%  '
$ws.Range("D46").Value = '2.068.20'
$ws.Range("E46").Value = '  +2.35%  '
$ws.Range("E48").Value = '  +8.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("E50").Value = '  +4.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '81.74'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.01%  '
